$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-30 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-31 Monday", 2) | Out-Null
$d.Content.Find.Execute("646÷4=161, 2", $true, $false, $false, $false, $false, $true, 1, $false, "892÷2=446, 0", 2) | Out-Null
$d.Content.Find.Execute("525÷8=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "825÷9=91, 6", 2) | Out-Null
$d.Content.Find.Execute("195÷6=32, 3", $true, $false, $false, $false, $false, $true, 1, $false, "798÷7=114, 0", 2) | Out-Null
$d.Content.Find.Execute("570÷7=81, 3", $true, $false, $false, $false, $false, $true, 1, $false, "994÷6=165, 4", 2) | Out-Null
$d.Content.Find.Execute("747÷3=249, 0", $true, $false, $false, $false, $false, $true, 1, $false, "966÷8=120, 6", 2) | Out-Null
$d.Content.Find.Execute("613÷2=306, 1", $true, $false, $false, $false, $false, $true, 1, $false, "284÷8=35, 4", 2) | Out-Null
$d.Content.Find.Execute("394÷2=197, 0", $true, $false, $false, $false, $false, $true, 1, $false, "382÷2=191, 0", 2) | Out-Null
$d.Content.Find.Execute("704÷9=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "823÷7=117, 4", 2) | Out-Null
$d.Content.Find.Execute("233÷2=116, 1", $true, $false, $false, $false, $false, $true, 1, $false, "291÷2=145, 1", 2) | Out-Null
$d.Content.Find.Execute("413÷2=206, 1", $true, $false, $false, $false, $false, $true, 1, $false, "548÷6=91, 2", 2) | Out-Null
$d.Content.Find.Execute("238÷6=39, 4", $true, $false, $false, $false, $false, $true, 1, $false, "484÷5=96, 4", 2) | Out-Null
$d.Content.Find.Execute("993÷8=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "101÷8=12, 5", 2) | Out-Null
$d.Content.Find.Execute("217÷2=108, 1", $true, $false, $false, $false, $false, $true, 1, $false, "858÷2=429, 0", 2) | Out-Null
$d.Content.Find.Execute("859÷6=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "374÷6=62, 2", 2) | Out-Null
$d.Content.Find.Execute("436÷5=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "567÷9=63, 0", 2) | Out-Null
$d.Content.Find.Execute("898÷4=224, 2", $true, $false, $false, $false, $false, $true, 1, $false, "100÷5=20, 0", 2) | Out-Null
$d.Content.Find.Execute("627÷9=69, 6", $true, $false, $false, $false, $false, $true, 1, $false, "866÷9=96, 2", 2) | Out-Null
$d.Content.Find.Execute("989÷8=123, 5", $true, $false, $false, $false, $false, $true, 1, $false, "520÷8=65, 0", 2) | Out-Null
$d.Content.Find.Execute("117÷2=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "415÷8=51, 7", 2) | Out-Null
$d.Content.Find.Execute("109÷4=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "752÷5=150, 2", 2) | Out-Null
$d.Content.Find.Execute("105÷9=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "362÷9=40, 2", 2) | Out-Null
$d.Content.Find.Execute("432÷4=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "449÷3=149, 2", 2) | Out-Null
$d.Content.Find.Execute("919÷9=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "505÷8=63, 1", 2) | Out-Null
$d.Content.Find.Execute("847÷8=105, 7", $true, $false, $false, $false, $false, $true, 1, $false, "442÷5=88, 2", 2) | Out-Null
$d.Content.Find.Execute("758÷6=126, 2", $true, $false, $false, $false, $false, $true, 1, $false, "388÷9=43, 1", 2) | Out-Null
